# Adjust stochastic example to actually be different for Scenario A & B
# Only ScenarioB's VRES data is changed so the two scenarios differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioB")

$ws.Range("F8").Value = 10
$ws.Range("G9").Value = 70
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 24
$ws.Range("G11").Value = 130
